# Auto-generated PowerShell COM-interop script
# Applies the market-price refresh diff to the Sheets workbook (Ramuh_Profits).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3948.843
$ws.Range("I69").Value = 4106.25
$ws.Range("J69").Value = 3876.8857
$ws.Range("K69").Value = 12318.75
$ws.Range("L69").Value = 11630.6571
$ws.Range("M69").Value = -11444.75
$ws.Range("N69").Value = -13378.6571
$ws.Range("H72").Value = 3948.843
$ws.Range("I72").Value = 4106.25
$ws.Range("J72").Value = 3876.8857
$ws.Range("K72").Value = 36956.25
$ws.Range("L72").Value = 34891.9713
$ws.Range("M72").Value = -32588.25
$ws.Range("N72").Value = -43627.9713
$ws.Range("H121").Value = 1508.75
$ws.Range("I121").Value = 400
$ws.Range("J121").Value = 1730.5
$ws.Range("K121").Value = 1200
$ws.Range("L121").Value = 5191.5
$ws.Range("M121").Value = 547
$ws.Range("N121").Value = -8685.5
$ws.Range("H137").Value = 4475.4595
$ws.Range("I137").Value = 1267.9231
$ws.Range("J137").Value = 6212.875
$ws.Range("K137").Value = 3803.7693
$ws.Range("L137").Value = 18638.625
$ws.Range("M137").Value = -1253.7693
$ws.Range("N137").Value = -23738.625

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 16951492
$ws.Range("I61").Value = 21741338
$ws.Range("K61").Value = 21741338
$ws.Range("M61").Value = -21741126
$ws.Range("H74").Value = 9185.154
$ws.Range("I74").Value = 1008.6667
$ws.Range("J74").Value = 16193.571
$ws.Range("K74").Value = 1008.6667
$ws.Range("L74").Value = 16193.571
$ws.Range("M74").Value = -134.6667
$ws.Range("N74").Value = -17941.571
$ws.Range("H77").Value = 9185.154
$ws.Range("I77").Value = 1008.6667
$ws.Range("J77").Value = 16193.571
$ws.Range("K77").Value = 5043.3335
$ws.Range("L77").Value = 80967.855
$ws.Range("M77").Value = -675.3334999999997
$ws.Range("N77").Value = -89703.855
$ws.Range("H132").Value = 5822.0713
$ws.Range("I132").Value = 1620.5883
$ws.Range("J132").Value = 12315.272
$ws.Range("K132").Value = 4861.7649
$ws.Range("L132").Value = 36945.81600000001
$ws.Range("M132").Value = -2331.7649
$ws.Range("N132").Value = -42005.81600000001
$ws.Range("H136").Value = 16951492
$ws.Range("I136").Value = 21741338
$ws.Range("K136").Value = 65224014
$ws.Range("M136").Value = -65221464

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1724.3334
$ws.Range("I134").Value = 1687.4286
$ws.Range("K134").Value = 5062.2858
$ws.Range("M134").Value = -2527.2858

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28605252
$ws.Range("I31").Value = 66668164
$ws.Range("J31").Value = 58067.8
$ws.Range("K31").Value = 66668164
$ws.Range("L31").Value = 58067.8
$ws.Range("M31").Value = -66667869
$ws.Range("N31").Value = -58657.8
$ws.Range("H34").Value = 28605252
$ws.Range("I34").Value = 66668164
$ws.Range("J34").Value = 58067.8
$ws.Range("K34").Value = 66668164
$ws.Range("L34").Value = 58067.8
$ws.Range("M34").Value = -66667962
$ws.Range("N34").Value = -58471.8
$ws.Range("H58").Value = 1509.9445
$ws.Range("I58").Value = 1585.7693
$ws.Range("J58").Value = 1312.8
$ws.Range("K58").Value = 1585.7693
$ws.Range("L58").Value = 1312.8
$ws.Range("M58").Value = -1382.7693
$ws.Range("N58").Value = -1718.8
$ws.Range("H132").Value = 25642948
$ws.Range("I132").Value = 33334766
$ws.Range("J132").Value = 3553.111
$ws.Range("K132").Value = 100004298
$ws.Range("L132").Value = 10659.333
$ws.Range("M132").Value = -100001768
$ws.Range("N132").Value = -15719.333
$ws.Range("H134").Value = 2623.913
$ws.Range("I134").Value = 2628.1333
$ws.Range("K134").Value = 7884.3999
$ws.Range("M134").Value = -5349.3999
$ws.Range("H136").Value = 1509.9445
$ws.Range("I136").Value = 1585.7693
$ws.Range("J136").Value = 1312.8
$ws.Range("K136").Value = 4757.3079
$ws.Range("L136").Value = 3938.4
$ws.Range("M136").Value = -2207.3079
$ws.Range("N136").Value = -9038.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 6303.7104
$ws.Range("I121").Value = 364
$ws.Range("J121").Value = 6633.6943
$ws.Range("K121").Value = 1092
$ws.Range("L121").Value = 19901.0829
$ws.Range("M121").Value = 218
$ws.Range("N121").Value = -22521.0829

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2984.0527
$ws.Range("I80").Value = 2302.8572
$ws.Range("J80").Value = 3381.4167
$ws.Range("K80").Value = 2302.8572
$ws.Range("L80").Value = 3381.4167
$ws.Range("M80").Value = -1304.8572
$ws.Range("N80").Value = -5377.4167
$ws.Range("H83").Value = 2984.0527
$ws.Range("I83").Value = 2302.8572
$ws.Range("J83").Value = 3381.4167
$ws.Range("K83").Value = 11514.286
$ws.Range("L83").Value = 16907.0835
$ws.Range("M83").Value = -6522.286
$ws.Range("N83").Value = -26891.0835
$ws.Range("H132").Value = 2504.5483
$ws.Range("I132").Value = 1972.2
$ws.Range("J132").Value = 3472.4546
$ws.Range("K132").Value = 5916.6
$ws.Range("L132").Value = 10417.3638
$ws.Range("M132").Value = -3386.6
$ws.Range("N132").Value = -15477.3638

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1392.5
$ws.Range("I22").Value = 1165.6154
$ws.Range("J22").Value = 1813.8572
$ws.Range("K22").Value = 1165.6154
$ws.Range("L22").Value = 1813.8572
$ws.Range("M22").Value = -870.6153999999999
$ws.Range("N22").Value = -2403.8572
$ws.Range("H27").Value = 1392.5
$ws.Range("I27").Value = 1165.6154
$ws.Range("J27").Value = 1813.8572
$ws.Range("K27").Value = 1165.6154
$ws.Range("L27").Value = 1813.8572
$ws.Range("M27").Value = -1058.6154
$ws.Range("N27").Value = -2027.8572
$ws.Range("H68").Value = 1730.3334
$ws.Range("J68").Value = 2145
$ws.Range("L68").Value = 2145
$ws.Range("N68").Value = -3643
$ws.Range("H71").Value = 1730.3334
$ws.Range("J71").Value = 2145
$ws.Range("L71").Value = 10725
$ws.Range("N71").Value = -18213
$ws.Range("H82").Value = 1750
$ws.Range("J82").Value = 1750
$ws.Range("L82").Value = 1750
$ws.Range("N82").Value = -2472
$ws.Range("H85").Value = 1750
$ws.Range("J85").Value = 1750
$ws.Range("L85").Value = 1750
$ws.Range("N85").Value = -4246
$ws.Range("H97").Value = 19036
$ws.Range("J97").Value = 19036
$ws.Range("L97").Value = 19036
$ws.Range("N97").Value = -21018
$ws.Range("H132").Value = 3064.8125
$ws.Range("I132").Value = 1866.6666
$ws.Range("J132").Value = 3341.3076
$ws.Range("K132").Value = 5599.9998
$ws.Range("L132").Value = 10023.9228
$ws.Range("M132").Value = -3069.9998
$ws.Range("N132").Value = -15083.9228

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2730.0454
$ws.Range("I62").Value = 2769.1428
$ws.Range("K62").Value = 2769.1428
$ws.Range("M62").Value = -2145.1428
$ws.Range("H65").Value = 2730.0454
$ws.Range("I65").Value = 2769.1428
$ws.Range("K65").Value = 13845.714
$ws.Range("M65").Value = -10725.714
$ws.Range("H94").Value = 34000
$ws.Range("J94").Value = 34000
$ws.Range("L94").Value = 34000
$ws.Range("N94").Value = -35802
$ws.Range("H95").Value = 26333.334
$ws.Range("J95").Value = 26333.334
$ws.Range("L95").Value = 26333.334
$ws.Range("N95").Value = -31825.334
$ws.Range("H122").Value = 11112605
$ws.Range("I122").Value = 18183104
$ws.Range("J122").Value = 1820
$ws.Range("K122").Value = 54549312
$ws.Range("L122").Value = 5460
$ws.Range("M122").Value = -54546862
$ws.Range("N122").Value = -10360
$ws.Range("H132").Value = 20002056
$ws.Range("I132").Value = 29413576
$ws.Range("J132").Value = 2575.375
$ws.Range("K132").Value = 88240728
$ws.Range("L132").Value = 7726.125
$ws.Range("M132").Value = -88238198
$ws.Range("N132").Value = -12786.125
$ws.Range("H136").Value = 5681.64
$ws.Range("I136").Value = 8686.615
$ws.Range("J136").Value = 2426.25
$ws.Range("K136").Value = 26059.845
$ws.Range("L136").Value = 7278.75
$ws.Range("M136").Value = -23509.845
$ws.Range("N136").Value = -12378.75
